$d = $word.ActiveDocument

# The template's address cell held the placeholder street number "333"
# (used twice in this photo-request template). Update it to the real
# address "476 loja 4" everywhere it occurs.
$d.Content.Find.Execute("333", $true, $true, $false, $false, $false, `
                         $true, 1, $false, "476 loja 4", 2) | Out-Null

$d.Save()
